$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Archive the current ("old") inference results (column M, rows 3-19) and
#    the old RMSE summary block (P23:R24) further down the sheet before they
#    get overwritten with the freshly recomputed values. This mirrors what
#    happened in the workbook: the previous results were copied down as a
#    reference/backup, then the M column was refreshed with the corrected
#    model output.
# ---------------------------------------------------------------------------

# Old M column values (rows 3-19), preserved at rows 22-38 (offset +19)
$oldM = @{
    3  = 55.680476331125107
    4  = 41.756478138640432
    5  = 84.929670372295817
    6  = 37.825875345796227
    7  = 35.751889036750981
    8  = 34.905428976361989
    9  = 33.078834354332727
    10 = 173.98223948728511
    11 = 91.539508452843435
    12 = 59.626571692041168
    13 = 67.174074696381723
    14 = 64.088337109654304
    15 = 109.9056922983319
    16 = 473.38039538904297
    17 = 420.38506339881877
    18 = 41.900994858328531
    19 = 39.768322133434268
}

# Rows whose old M cell carried the shaded "s=3" style; that formatting is
# carried over to the archived copy.
$styledRows = @(4, 5, 6, 12, 14, 15, 16, 17)

for ($r = 3; $r -le 19; $r++) {
    $destRow = $r + 19
    $ws.Cells.Item($destRow, 13).Value = $oldM[$r]
    if ($styledRows -contains $r) {
        $ws.Cells.Item($r, 13).Copy()
        $ws.Cells.Item($destRow, 13).PasteSpecial(-4122)
        $ws.Cells.Item($destRow, 13).Value = $oldM[$r]
    }
}
$excel.CutCopyMode = 0

# Archive the old RMSE comparison block (P23:R24) as plain values two rows
# below its previous location (P31:R32), without formulas or formatting.
$ws.Cells.Item(31, 16).Value = $ws.Cells.Item(23, 16).Value()
$ws.Cells.Item(31, 17).Value = $ws.Cells.Item(23, 17).Value()
$ws.Cells.Item(31, 18).Value = $ws.Cells.Item(23, 18).Value()
$ws.Cells.Item(32, 16).Value = $ws.Cells.Item(24, 16).Value()
$ws.Cells.Item(32, 17).Value = $ws.Cells.Item(24, 17).Value()
$ws.Cells.Item(32, 18).Value = $ws.Cells.Item(24, 18).Value()

# ---------------------------------------------------------------------------
# 2. Write the corrected/recomputed inference results into column M.
#    Several of these cells also lose their old shaded "s=3" formatting
#    (the correction pass cleared it for rows 4,5,6,12,14,15,16,17).
# ---------------------------------------------------------------------------
$newM = @{
    3  = 54.605536396961817
    4  = 41.856676411395242
    5  = 87.577477027577629
    6  = 37.377088243648622
    7  = 35.401481552743753
    8  = 34.467716296422438
    9  = 32.418480153017569
    10 = 179.74905279769121
    11 = 94.509138300523929
    12 = 58.623068603283421
    13 = 67.74274199969517
    14 = 64.391840727219488
    15 = 110.8674542020869
    16 = 499.33953049095072
    17 = 442.80369324368343
    18 = 42.046688141426777
    19 = 38.934461128219553
}

foreach ($r in $styledRows) {
    $ws.Cells.Item($r, 13).ClearFormats()
}

for ($r = 3; $r -le 19; $r++) {
    $ws.Cells.Item($r, 13).Value = $newM[$r]
}

# ---------------------------------------------------------------------------
# 3. Restyle N19 to match the rest of row 19 (bottom border, no fill),
#    copying the existing format from B19.
# ---------------------------------------------------------------------------
$ws.Cells.Item(19, 2).Copy()
$ws.Cells.Item(19, 14).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Update the active selection to reflect where the edit left off.
# ---------------------------------------------------------------------------
$ws.Range("K23").Select() | Out-Null
